$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal TEXT value into a cell even when the text looks
# like a plain number (e.g. "0.17", "-0.89"), without Excel's normal
# "looks-like-a-number -> store as number" auto-conversion, and without
# leaving behind a changed NumberFormat/style on the cell (matches the
# original file, where every one of these values is a shared string with
# the sheet's default style).
#
# Trick: put a formula that evaluates to the text string in a scratch
# cell, copy it, and paste-special as values into the destination. A
# formula result of type string pastes as a literal text cell (t="str"
# becomes t="s" on save) instead of being re-parsed as a number.
function Set-TextValue($addr, $val) {
    $scratch = $ws.Range("Z100")
    $escaped = $val.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Value = $null
}

# Write column-by-column (B then C then D) to keep the shared-string table
# ordered the same way the source file had it (A Lag/FFR Lag/LF Lag rows
# within each column, column by column).
Set-TextValue "B2" "0.17"
Set-TextValue "B3" "-0.01"
Set-TextValue "B4" "-0.09"

$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"
Set-TextValue "C4" "0.98"

Set-TextValue "D2" "-0.89"
$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"

$excel.CutCopyMode = $false
